$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = [double]"2.475797344914099e-14"
$ws.Range("I2").Value = [double]"2.475797344914099e-14"
$ws.Range("L2").Value = [double]"58.54608633267811"
$ws.Range("M2").Value = "[47.36139855250467, 69.73077411285155]"
$ws.Range("N2").Value = [double]"9.681144774731365e-14"
$ws.Range("O2").Value = [double]"9.681144774731365e-14"
$ws.Range("P2").Value = [double]"1.440289725069195"
$ws.Range("Q2").Value = "[1.2264475824825025, 1.6541318676558872]"
$ws.Range("T2").Value = [double]"49.98199184545564"
$ws.Range("U2").Value = "[42.9592257927675, 57.00475789814378]"
$ws.Range("V2").Value = [double]"0"
$ws.Range("W2").Value = [double]"0"
$ws.Range("X2").Value = [double]"17.8587587587589"
$ws.Range("Y2").Value = [double]"17.07019019019032"
$ws.Range("Z2").Value = [double]"18.64732732732747"
$ws.Range("H3").Value = [double]"2.783551167340192e-12"
$ws.Range("I3").Value = [double]"2.783551167340192e-12"
$ws.Range("L3").Value = [double]"53.21247739906848"
$ws.Range("M3").Value = "[40.48556725333137, 65.93938754480558]"
$ws.Range("N3").Value = [double]"8.472000878612107e-11"
$ws.Range("O3").Value = [double]"8.472000878612107e-11"
$ws.Range("P3").Value = [double]"1.767342413731195"
$ws.Range("Q3").Value = "[1.515763422452733, 2.0189214050096576]"
$ws.Range("T3").Value = [double]"50.0009716250046"
$ws.Range("U3").Value = "[42.604293628370705, 57.3976496216385]"
$ws.Range("X3").Value = [double]"16.65271271271284"
$ws.Range("Y3").Value = [double]"15.72498498498511"
$ws.Range("Z3").Value = [double]"17.58044044044058"
$ws.Range("H4").Value = [double]"2.919464670014804e-11"
$ws.Range("I4").Value = [double]"2.919464670014804e-11"
$ws.Range("L4").Value = [double]"63.17791653770873"
$ws.Range("M4").Value = "[45.18965184265281, 81.16618123276464]"
$ws.Range("N4").Value = [double]"7.918861566480473e-09"
$ws.Range("O4").Value = [double]"7.918861566480473e-09"
$ws.Range("P4").Value = [double]"2.094395102393196"
$ws.Range("Q4").Value = "[1.7925003128590413, 2.39628989192735]"
$ws.Range("T4").Value = [double]"57.70387898545621"
$ws.Range("U4").Value = "[48.145393185349064, 67.26236478556336]"
$ws.Range("V4").Value = [double]"8.881784197001252e-16"
$ws.Range("W4").Value = [double]"8.881784197001252e-16"
$ws.Range("X4").Value = [double]"15.44666666666679"
$ws.Range("Y4").Value = [double]"14.3333933933935"
$ws.Range("Z4").Value = [double]"16.55993993994007"
$ws.Range("H5").Value = [double]"1.975675179011205e-11"
$ws.Range("I5").Value = [double]"1.975675179011205e-11"
$ws.Range("L5").Value = [double]"52.30573581471734"
$ws.Range("M5").Value = "[37.35256207211066, 67.25890955732403]"
$ws.Range("N5").Value = [double]"8.730876910334473e-09"
$ws.Range("O5").Value = [double]"8.730876910334473e-09"
$ws.Range("P5").Value = [double]"2.169868799776734"
$ws.Range("Q5").Value = "[1.8679740102425795, 2.4717635893108882]"
$ws.Range("T5").Value = [double]"47.3586615359648"
$ws.Range("U5").Value = "[39.521117840306005, 55.196205231623594]"
$ws.Range("V5").Value = [double]"6.661338147750939e-16"
$ws.Range("W5").Value = [double]"6.661338147750939e-16"
$ws.Range("X5").Value = [double]"15.16834834834847"
$ws.Range("Y5").Value = [double]"14.05507507507519"
$ws.Range("Z5").Value = [double]"16.28162162162175"
$ws.Range("H6").Value = [double]"1.432187701766452e-14"
$ws.Range("I6").Value = [double]"1.432187701766452e-14"
$ws.Range("L6").Value = [double]"57.75147108771576"
$ws.Range("M6").Value = "[43.89428330635461, 71.6086588690769]"
$ws.Range("N6").Value = [double]"9.26958509950282e-11"
$ws.Range("O6").Value = [double]"9.26958509950282e-11"
$ws.Range("P6").Value = [double]"2.333395144107735"
$ws.Range("Q6").Value = "[2.1069740519571187, 2.5598162362583503]"
$ws.Range("T6").Value = [double]"53.35184280844373"
$ws.Range("U6").Value = "[46.41527078103088, 60.28841483585659]"
$ws.Range("X6").Value = [double]"14.56532532532544"
$ws.Range("Y6").Value = [double]"13.73037037037048"
$ws.Range("Z6").Value = [double]"15.4002802802804"
$ws.Range("H7").Value = [double]"4.430900091279e-13"
$ws.Range("I7").Value = [double]"4.430900091279e-13"
$ws.Range("L7").Value = [double]"60.32294167392398"
$ws.Range("M7").Value = "[45.04158624468138, 75.60429710316657]"
$ws.Range("N7").Value = [double]"4.067406411678576e-10"
$ws.Range("O7").Value = [double]"4.067406411678576e-10"
$ws.Range("P7").Value = [double]"2.635289933641888"
$ws.Range("Q7").Value = "[2.3711319927995036, 2.8994478744842733]"
$ws.Range("T7").Value = [double]"53.36600118331512"
$ws.Range("U7").Value = "[45.251722245264766, 61.480280121365475]"
$ws.Range("X7").Value = [double]"13.45205205205216"
$ws.Range("Y7").Value = [double]"12.47793793793804"
$ws.Range("Z7").Value = [double]"14.42616616616628"
$ws.Range("H8").Value = [double]"2.55351295663786e-14"
$ws.Range("I8").Value = [double]"2.55351295663786e-14"
$ws.Range("L8").Value = [double]"63.00032225104706"
$ws.Range("M8").Value = "[49.348616624195074, 76.65202787789904]"
$ws.Range("N8").Value = [double]"4.883426996116214e-12"
$ws.Range("O8").Value = [double]"4.883426996116214e-12"
$ws.Range("P8").Value = [double]"2.899447874484274"
$ws.Range("Q8").Value = "[2.660447832769735, 3.1384479161988126]"
$ws.Range("T8").Value = [double]"55.92303910862606"
$ws.Range("U8").Value = "[48.105649477093586, 63.74042874015853]"
$ws.Range("X8").Value = [double]"12.47793793793804"
$ws.Range("Y8").Value = [double]"11.59659659659669"
$ws.Range("Z8").Value = [double]"13.35927927927938"
$ws.Range("F9").Value = [double]"22.6300000000001"
$ws.Range("H9").Value = [double]"1.684652417566213e-11"
$ws.Range("I9").Value = [double]"1.684652417566213e-11"
$ws.Range("L9").Value = [double]"52.71228075331926"
$ws.Range("M9").Value = "[39.86176628813348, 65.56279521850503]"
$ws.Range("N9").Value = [double]"1.43817624476128e-10"
$ws.Range("O9").Value = [double]"1.43817624476128e-10"
$ws.Range("P9").Value = [double]"3.037816319687427"
$ws.Range("Q9").Value = "[2.76107942928112, 3.314553210093735]"
$ws.Range("T9").Value = [double]"52.07592316978028"
$ws.Range("U9").Value = "[43.97151395652623, 60.18033238303433]"
$ws.Range("V9").Value = [double]"2.220446049250313e-16"
$ws.Range("W9").Value = [double]"2.220446049250313e-16"
$ws.Range("X9").Value = [double]"11.68876876876882"
$ws.Range("Y9").Value = [double]"10.6920520520521"
$ws.Range("Z9").Value = [double]"12.68548548548554"
$ws.Range("F10").Value = [double]"22.6300000000001"
$ws.Range("H10").Value = [double]"1.998401444325282e-15"
$ws.Range("I10").Value = [double]"1.998401444325282e-15"
$ws.Range("L10").Value = [double]"58.30885175301719"
$ws.Range("M10").Value = "[47.3129774436797, 69.30472606235469]"
$ws.Range("N10").Value = [double]"6.350475700855895e-14"
$ws.Range("O10").Value = [double]"6.350475700855895e-14"
$ws.Range("P10").Value = [double]"-3.00636894577762"
$ws.Range("Q10").Value = "[-3.20763213880039, -2.8051057527548497]"
$ws.Range("T10").Value = [double]"53.81026391962709"
$ws.Range("U10").Value = "[47.03484978698536, 60.58567805226882]"
$ws.Range("X10").Value = [double]"10.82796796796801"
$ws.Range("Y10").Value = [double]"10.10308308308312"
$ws.Range("Z10").Value = [double]"11.5528528528529"
$ws.Range("F11").Value = [double]"22.6300000000001"
$ws.Range("H11").Value = [double]"4.561906408184768e-13"
$ws.Range("I11").Value = [double]"4.561906408184768e-13"
$ws.Range("L11").Value = [double]"57.40193008241496"
$ws.Range("M11").Value = "[42.31916078561868, 72.48469937921124]"
$ws.Range("N11").Value = [double]"1.063010124369157e-09"
$ws.Range("O11").Value = [double]"1.063010124369157e-09"
$ws.Range("P11").Value = [double]"-2.490632013656773"
$ws.Range("Q11").Value = "[-2.754789954499157, -2.2264740728143884]"
$ws.Range("T11").Value = [double]"54.87453515276783"
$ws.Range("U11").Value = "[47.20775455247836, 62.541315753057305]"
$ws.Range("X11").Value = [double]"8.970450450450491"
$ws.Range("Y11").Value = [double]"8.019039039039077"
$ws.Range("Z11").Value = [double]"9.921861861861904"
$ws.Range("F12").Value = [double]"22.6300000000001"
$ws.Range("H12").Value = [double]"1.865174681370263e-13"
$ws.Range("I12").Value = [double]"1.865174681370263e-13"
$ws.Range("L12").Value = [double]"60.32216082881818"
$ws.Range("M12").Value = "[45.665959232130746, 74.97836242550562]"
$ws.Range("N12").Value = [double]"1.310715980196164e-10"
$ws.Range("O12").Value = [double]"1.310715980196164e-10"
$ws.Range("P12").Value = [double]"-2.163579324994773"
$ws.Range("Q12").Value = "[-2.402579366709311, -1.9245792832802344]"
$ws.Range("R12").Value = [double]"0"
$ws.Range("S12").Value = [double]"0"
$ws.Range("T12").Value = [double]"50.55037092662859"
$ws.Range("U12").Value = "[42.83861263310985, 58.26212922014733]"
$ws.Range("V12").Value = [double]"0"
$ws.Range("W12").Value = [double]"0"
$ws.Range("X12").Value = [double]"7.792512512512548"
$ws.Range("Y12").Value = [double]"6.931711711711745"
$ws.Range("Z12").Value = [double]"8.653313313313351"
